$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = -0.000052956794893893
$ws.Range("B3").Value = 0.7929965008397935
$ws.Range("B4").Value = 0.001978262929895643
$ws.Range("B5").Value = "[1.0, 0.31967756971198075, 0.3038434045935043, 0.2941671851070525, 0.2727377628275661, 0.47359109975194336, 0.27698584982002694, 0.294941860423538, 0.2825140517732761, 0.2467495898746261, 0.23526531121076355, 0.22260491191002388, 0.24758520241450924, 0.24618261009451964, 0.2623274664137683, 0.2393553527768477, 0.20773651918313554, 0.23385849195393235, 0.21116232515339256, 0.2057272454195783]"

$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B3").Value = 1.067583974099555
$ws.Range("B4").Value = -0.01163059438434978
$ws.Range("B5").Value = 0.7517595727673099
$ws.Range("B6").Value = 0.001207455341984045
$ws.Range("B7").Value = 0.01256828842261183
$ws.Range("B8").Value = 0.7462214371342081
$ws.Range("B9").Value = 0.001444694635368234
$ws.Range("B10").Value = "[0.9999999999999999, 0.3457942028043245, 0.30187447422524244, 0.32874727556479366, 0.28005497474281904, 0.34290257139185665, 0.2710303012689809, 0.2722019259068867, 0.2692478390950436, 0.23395127770014754, 0.24595626206923216, 0.20977900467558047, 0.22395637177089922, 0.20446574399531, 0.23136328413254656, 0.24503618823815235, 0.2025509470082138, 0.21427410187973012, 0.20758283919148882, 0.19333797985146114]"

$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = 0.0003351416656164028
$ws.Range("B3").Value = 0.00004413315177652508
$ws.Range("B4").Value = 0.1999999745631566
$ws.Range("B5").Value = 0.7800000285684567
$ws.Range("B6").Value = "[0.9999999999999999, -0.09021303804169639, -0.05069690986088107, -0.07357237745612387, -0.11834949268713568, 0.28775386407957737, -0.04668057576844395, 0.05852392703767862, 0.0036647847892249106, -0.023474008545813524, -0.038940786472992366, -0.03772287650368242, 0.016993802621373847, 0.021073064777259242, 0.06962543253462977, 0.012218773106454854, -0.003873741031019354, -0.010859006176298293, -0.028785102078740535, 0.004470907994055973]"

$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.0003438700948161704
$ws.Range("B3").Value = 0.00004413636892290445
$ws.Range("B4").Value = 0.1000000889766707
$ws.Range("B5").Value = 0.8549999255094036
$ws.Range("B6").Value = "[1.0, -0.04537073463150119, -0.014809631493741164, -0.039996503971106684, -0.08729354845500342, 0.30487741041294214, -0.02139768261892874, 0.07460534623267191, 0.022457070909931, -0.013671819562525198, -0.027103267391320216, -0.027760763788982637, 0.026844212946170516, 0.03078091553296756, 0.06966426059382526, 0.019116851364118073, 0.0009243541293403281, -0.004575229639319514, -0.022159726163030472, 0.009127335732338376]"
$ws.Range("B7").Value = 0.04999986672740196

$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = -0.0003176544666566077
$ws.Range("B3").Value = 0.00003956687539925694
$ws.Range("B4").Value = 0.199983684524659
$ws.Range("B5").Value = 0.7750186527439975
$ws.Range("B6").Value = "[1.0, -0.0676488815431104, -0.031416736919497305, -0.047700995538358874, -0.06530600588951523, 0.1998446210289016, -0.03509617009270588, 0.03344754199212242, -0.005953713292294641, -0.021100965572721096, -0.02579211762132092, -0.018559131795802725, 0.009000413341306314, 0.003330477239898005, 0.055756450823325344, 0.016599528182636426, -0.006120883524487295, 0.006796309748721854, -0.002617752846073729, 0.0024843077130455374]"
$ws.Range("B7").Value = 0.009999222019920005
$ws.Range("B9").Value = 0.7604946193121846
